$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = 0.444444
$ws.Range("B69").Value = 0.333333
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 0.888888
$ws.Range("E69").Value = 0.999999
$ws.Range("F69").Value = 0.111111
$ws.Range("G69").Value = 0.08538695911315511
$ws.Range("H69").Value = "query"

$ws.Range("A70").Value = 0.222222
$ws.Range("B70").Value = 0.222222
$ws.Range("C70").Value = 0.7777770000000001
$ws.Range("D70").Value = 0.333333
$ws.Range("E70").Value = 0.444444
$ws.Range("F70").Value = 0.111111
$ws.Range("G70").Value = 0.1937091356989436
$ws.Range("H70").Value = "query"
